$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to force text-typed entry for values that would
# otherwise be auto-coerced to numbers by Excel (e.g. "321.84").
# NumberFormat "@" on the scratch cell marks the typed value as text;
# PasteSpecial(xlPasteValues) then carries only the value (not the
# format) into the destination cell, leaving its style untouched.
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

function Set-TextValue([string]$addr, [string]$text) {
    $scratch.Value = $text
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163) # xlPasteValues
}

$ws.Range("D2").Value = '48.283.29'
$ws.Range("E2").Value = '  +1.55%  '
$ws.Range("D3").Value = '2.508.36'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("E4").Value = '  +0.01%  '
Set-TextValue "D5" '321.84'
$ws.Range("E5").Value = '  -0.03%  '
Set-TextValue "D6" '108.13'
$ws.Range("E6").Value = '  -0.82%  '
$ws.Range("E7").Value = '  +1.12%  '
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -0.06%  '
Set-TextValue "D10" '39.83'
$ws.Range("E10").Value = '  +1.43%  '
Set-TextValue "D11" '20.26'
$ws.Range("E11").Value = '  +8.68%  '
Set-TextValue "D12" '0.0820'
$ws.Range("E12").Value = '  +1.17%  '
$ws.Range("E13").Value = '  +0.06%  '
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("D16").Value = '2.508.29'
$ws.Range("E16").Value = '  +0.38%  '
Set-TextValue "D17" '0.845'
$ws.Range("E17").Value = '  -0.48%  '
$ws.Range("D18").Value = '48.143.16'
$ws.Range("E18").Value = '  +1.48%  '
Set-TextValue "D19" '13.12'
$ws.Range("E19").Value = '  -1.83%  '
Set-TextValue "D20" '6.80'
$ws.Range("E20").Value = '  +2.30%  '
$ws.Range("D21").Value = '0.0₃0945'
Set-TextValue "D22" '2.77'
$ws.Range("E22").Value = '  +0.52%  '
Set-TextValue "D23" '280.74'
$ws.Range("E23").Value = '  +13.82%  '
Set-TextValue "D24" '72.32'
$ws.Range("E24").Value = '  +2.35%  '
$ws.Range("E25").Value = '  +0.01%  '
$ws.Range("E26").Value = '  -0.10%  '
Set-TextValue "D27" '25.77'
$ws.Range("E27").Value = '  +0.09%  '
Set-TextValue "D28" '2.21'
$ws.Range("E28").Value = '  +0.41%  '
$ws.Range("E29").Value = '  -1.99%  '
$ws.Range("E30").Value = '  +0.87%  '
Set-TextValue "D31" '35.32'
$ws.Range("E31").Value = '  +1.71%  '
Set-TextValue "D32" '49.51'
$ws.Range("E32").Value = '  -0.84%  '
Set-TextValue "D33" '19.72'
$ws.Range("E33").Value = '  -2.80%  '
$ws.Range("E34").Value = '  +0.26%  '
Set-TextValue "D35" '1.01'
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  -0.61%  '
$ws.Range("E37").Value = '  -0.36%  '
$ws.Range("E38").Value = '  -1.78%  '
Set-TextValue "D39" '2.92'
$ws.Range("E39").Value = '  -1.21%  '
$ws.Range("E40").Value = '  -0.01%  '
Set-TextValue "D41" '122.52'
$ws.Range("E41").Value = '  +3.01%  '
$ws.Range("E42").Value = '  -0.05%  '
Set-TextValue "D43" '21.52'
$ws.Range("E43").Value = '  -4.22%  '
$ws.Range("E44").Value = '  +1.80%  '
$ws.Range("D45").Value = '2.016.98'
$ws.Range("E45").Value = '  +0.89%  '
Set-TextValue "D46" '3.18'
$ws.Range("E46").Value = '  +4.68%  '
$ws.Range("E47").Value = '  +3.26%  '
$ws.Range("E48").Value = '  -2.41%  '
Set-TextValue "D49" '9.04'
$ws.Range("E49").Value = '  -0.09%  '
$ws.Range("E50").Value = '  -0.89%  '
Set-TextValue "D51" '80.66'
$ws.Range("E51").Value = '  +3.99%  '

# Clean up the scratch cell so it leaves no trace in the saved sheet.
$scratch.Clear()
$excel.CutCopyMode = $false

